# Generate Report for Handback
# Update the timestamp values recorded in the handback status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first row.
$wsOverview.Range("G2").Value = "2016-08-22 15:15:46"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for the first row.
$wsZhCn.Range("H2").Value = "2016-08-22 15:15:41"
$wsZhCn.Range("K2").Value = "2016-08-22 15:15:59"

# de-de sheet: "Correspond Handback DateTime" for the first row.
$wsDeDe.Range("K2").Value = "2016-08-22 15:16:22"
